$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Q0)
$ws.Range("B3").Value = 0.2128381982599638
$ws.Range("C3").Value = 0.5858536052976883
$ws.Range("D3").Value = 0.5319869414474643
$ws.Range("E3").Value = 0.7293743493210221
$ws.Range("F3").Value = 0.701291567186506
$ws.Range("G3").Value = 96

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1618673208653933
$ws.Range("C4").Value = 0.5051966086236355
$ws.Range("D4").Value = 0.340115955331918
$ws.Range("E4").Value = 0.5831946118851905
$ws.Range("F4").Value = 0.566472252636492
$ws.Range("G4").Value = 46
